$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data for row 13 (B13 and C13), matching the new shared string
$ws.Range("B13").Value = 4.5
$ws.Range("C13").Value = "preprocessing session 4"

# Add the SUM formula in D12 (summing B10:B13), written with a leading space like the source file
$ws.Range("D12").Formula = "= SUM(B10:B13)"

# Update the selected cell to D13, matching the final saved selection
$ws.Range("D13").Select()
